{"js": "// Update the salutation's name from \"Dr. Krystal\" to \"Dr. Carter\", and\n// change the first \"blue\" (the run describing the highlight colour used for\n// the audit note, color 00B0F0) to \"cyan\". The second \"blue\" (the\n// \"blue regular font\" mention, color 215E99) is left untouched.\n\nconst body = context.document.body;\n\n// 1) Dear Dr. Krystal, -> Dear Dr. Carter,\nconst nameHits = body.search(\"Dr. Krystal\", { matchCase: true });\nnameHits.load(\"text\");\nawait context.sync();\n\nif (nameHits.items.length > 0) {\n  nameHits.items[0].insertText(\"Dr. Carter\", \"Replace\");\n  await context.sync();\n}\n\n// 2) \"...marked in blue in the manuscript.\" -> \"...marked in cyan...\"\n//    Only the first exact-case \"blue\" occurrence (colour 00B0F0) changes;\n//    the later \"blue regular font\" (colour 215E99) must stay \"blue\".\nconst blueHits = body.search(\"blue\", { matchCase: true, matchWholeWord: true });\nblueHits.load(\"text\");\nawait context.sync();\n\nfor (const hit of blueHits.items) {\n  hit.font.load(\"color\");\n}\nawait context.sync();\n\nconst target = blueHits.items.find(\n  (hit) => (hit.font.color || \"\").toUpperCase() === \"#00B0F0\"\n) || blueHits.items[0];\n\nif (target) {\n  target.insertText(\"cyan\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the salutation's name from \"Dr. Krystal\" to \"Dr. Carter\", and\n# change the first \"blue\" (the run describing the highlight colour used for\n# the audit note, color 00B0F0) to \"cyan\". The second \"blue\" (the\n# \"blue regular font\" mention, color 215E99) is left untouched.\n\n$d = $word.ActiveDocument\n\n# 1) Dear Dr. Krystal, -> Dear Dr. Carter,\n$nameRange = $d.Content\n$nameFind = $nameRange.Find\n$nameFind.Text = \"Dr. Krystal\"\n$nameFind.MatchCase = $true\n$nameFind.MatchWholeWord = $false\nif ($nameFind.Execute()) {\n    $nameRange.Text = \"Dr. Carter\"\n}\n\n# 2) \"...marked in blue in the manuscript.\" -> \"...marked in cyan...\"\n#    Only the first exact-case, whole-word \"blue\" (colour 00B0F0) changes;\n#    the later \"blue regular font\" (colour 215E99) must stay \"blue\".\n$blueRange = $d.Content\n$blueFind = $blueRange.Find\n$blueFind.Text = \"blue\"\n$blueFind.MatchCase = $true\n$blueFind.MatchWholeWord = $true\n\nwhile ($blueFind.Execute()) {\n    if ($blueRange.Font.Color -eq 15773696) {\n        $blueRange.Text = \"cyan\"\n        break\n    }\n}\n"}
